$p = $ppt.ActivePresentation

# --- Slide 13: "Tightly Coupled Presentation Logic..." title placeholder ---
# Re-set the text through a transient two-paragraph round trip so the
# stored paragraph no longer carries a stray trailing endParaRPr run
# (mirrors how PowerPoint drops it once the paragraph has been touched).
$slide13 = $p.Slides.Item(13)

$shpA = $slide13.Shapes.Item(1)
$trA = $shpA.TextFrame.TextRange
$trA.Text = "X`rY"
$trA2 = $shpA.TextFrame.TextRange
$trA2.Text = "Tightly Coupled Presentation Logic: One Presenter Knows About Another"

# --- Slide 13: "Simple Org Chart - v1" title ---
$shpB = $slide13.Shapes.Item(2)
$trB = $shpB.TextFrame.TextRange
$trB.Text = "X`rY"
$trB2 = $shpB.TextFrame.TextRange
$trB2.Text = "Simple Org Chart – v1"

# --- Slide 20: merge "Orchestration" + " With An Application Controller And "
#     into a single run (keeps the first run's formatting). ---
$slide20 = $p.Slides.Item(20)
$shpC = $slide20.Shapes.Item(7)
$trC = $shpC.TextFrame.TextRange
$subC = $trC.Characters(1, 49)
$subC.Text = "Orchestration With An Application Controller And "
